$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 - this shifts existing rows 5..46 down to 6..47
# (the D-column date style carries over onto the new row from the Insert).
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly price observation.
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = Get-Date -Year 2022 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112037
$ws.Range("G5").Value = "Cebollín"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 180
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 6500
$ws.Range("M5").Value = 6222
$ws.Range("N5").Value = "$/paquete 36 unidades"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 173
$ws.Range("Q5").Value = 36
$ws.Range("R5").Value = "Hortaliza"
